
# "Added Arten to Variablen" -- slide 2 ("Variablen") content placeholder:
# the existing "Arten" bullet gets five new sub-bullets listing the
# variable types: integer, float, String, character, boolean.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame.TextRange

# Re-write the whole placeholder text, keeping the existing bullets
# ("Deklaration", "Zuweisung von Werten", "Arten") and appending the new
# sub-bullets, one per paragraph (CR = new paragraph in PowerPoint OM).
$tf.Text = "Deklaration`rZuweisung von Werten`rArten`rinteger`rfloat`rString`rcharacter`rboolean"

# Demote the five new paragraphs to the second outline level (lvl="1").
$tf.Characters(40, 7).IndentLevel = 2   # integer
$tf.Characters(48, 5).IndentLevel = 2   # float
$tf.Characters(54, 6).IndentLevel = 2   # String
$tf.Characters(61, 9).IndentLevel = 2   # character
$tf.Characters(71, 7).IndentLevel = 2   # boolean

# Split the run on the first letter for "integer", "float", "character"
# and "boolean" (mirrors the run layout produced by PowerPoint's
# autocorrect/spellcheck pass on first entry of these words).
$tf.Characters(40, 1).Text = "i"
$tf.Characters(48, 1).Text = "f"
$tf.Characters(61, 1).Text = "c"
$tf.Characters(71, 1).Text = "b"
